$wb = $excel.ActiveWorkbook

$wsIn = $wb.Worksheets.Item("ProductLoanInput")
$wsOut = $wb.Worksheets.Item("ProductLoanOutput")

# Update the product name (same new value on both the input and output sheets)
$newProductName = "4282-MS-EI-DB-SAR-REC-CTRFD-RNI-INT-FFC-SAR-FFROP-DAILY-FIFR-1-MD-TR-1-1st"
$wsIn.Range("B1").Value = $newProductName
$wsOut.Range("B1").Value = $newProductName

# Update the short name - was numeric 4282, now textual "428w"
$wsIn.Range("B2").Value = "428w"

# Move the selection on the input sheet off of B17 and onto B3 (test case
# inter-dependency removal), and make the output sheet the active tab.
$wsIn.Range("B3").Select() | Out-Null
$wsOut.Activate() | Out-Null
